# "graphe ia design actuel"
#
# Sheet1's classification data (package names + their inheritance/cohesion
# sample counts) is refreshed for the current (DM3) design: the old
# "Transport"-package rows are replaced by the new "Application"/"Travel.*"
# package breakdown, and the I/A ratio formulas (columns B/C) are
# re-derived for every row (C was F/G, now F/(G+F) to match B's shape).
# Rows 12-19 are brand-new packages that have no samples yet, so their
# formulas legitimately evaluate to #DIV/0! (no data to average).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 2-11: label + underlying D/E/F/G sample counts ---
# (B/C formulas recalc automatically from these; re-set explicitly below too)
$ws.Range("A2").Value = "Application"
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1

$ws.Range("A3").Value = "Sessions"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 7

$ws.Range("A4").Value = "Utils"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 8

$ws.Range("A5").Value = "Commands"
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 23

$ws.Range("A6").Value = "Travel"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 40

$ws.Range("A7").Value = "Réservation"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5

$ws.Range("A8").Value = "Paiment"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 2

$ws.Range("A9").Value = "UI"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4

$ws.Range("A10").Value = "Travel.Facilities"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3

$ws.Range("A11").Value = "Travel.Trip"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 4

# --- Rows 12-19: brand-new packages, label only (no samples yet) ---
$ws.Range("A12").Value = "Travel.Vehicule"
$ws.Range("A13").Value = "Travel.Factories"
$ws.Range("A14").Value = "Travel.Forms"
$ws.Range("A15").Value = "Travel.Itineries"
$ws.Range("A16").Value = "Travel.Sections"
$ws.Range("A17").Value = "Travel.Places"
$ws.Range("A18").Value = "Travel."
$ws.Range("A19").Value = "Travel.VechiculeModels"

# --- B/C ratio formulas for every data row ---
# B = D/(E+D)              (unchanged shape)
# C = F/(G+F)              (was F/G before this edit)
$ws.Range("B2").Formula = "=D2/(E2+D2)"
$ws.Range("C2").Formula = "=F2/(G2+F2)"
for ($r = 3; $r -le 19; $r++) {
    $ws.Range("B$r").Formula = "=D$r/(E$r+D$r)"
    $ws.Range("C$r").Formula = "=F$r/(G$r+F$r)"
}

# --- View bits: column A widened for the longer "Travel.*" labels,
#     active cell moved to D13, sheet zoomed to 100% ---
$ws.Range("A:A").ColumnWidth = 22.7109375
$ws.Range("D13").Select()
$excel.ActiveWindow.Zoom = 100
